$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 49 and 50 need to be swapped (content + formatting), and the label
# that ends up in row 49 ("pharmacological treatment") is renamed to
# "pharmacological intervention".
#
# Use a scratch row far outside the used range to stage row 49's original
# content/formatting while row 50's content/formatting is moved up.

$scratchRow = 1000

$row49 = $ws.Range("A49:W49")
$row50 = $ws.Range("A50:W50")
$scratch = $ws.Range("A" + $scratchRow + ":W" + $scratchRow)

# Stash current row 49 (CHEBI:52210 / pharmacological role ...) in scratch.
# (Clear the destination first -- Copy() only overwrites cells that are
# populated in the source range, it does not blank out leftover cells.)
$scratch.Clear()
$row49.Copy($scratch)

# Move current row 50 (GMHO:0000262 / pharmacological treatment ...) up into row 49.
$row49.Clear()
$row50.Copy($row49)

# Move the stashed original row 49 down into row 50.
$row50.Clear()
$scratch.Copy($row50)

# Clear the scratch row (values + formatting) so nothing extraneous is left behind.
$scratch.Clear()

# Rename the label now sitting in row 49.
$ws.Cells.Item(49, 2).Value = "pharmacological intervention"

# Separately, rename "psychological treatment" to "psychological intervention" in row 60.
$ws.Cells.Item(60, 2).Value = "psychological intervention"
